$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (A5:C5 = 1, 2, 3), matching the existing pattern
# used in rows 2-4.
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3

# Leave the selection where Excel would after typing across the row with Tab,
# i.e. one cell to the right of the last entry.
$null = $ws.Range("D5").Select()
